$d = $word.ActiveDocument

# The document currently opens with a lone paragraph, immediately
# before the first table, whose only content is a single stray "H"
# character (tiny 6pt Calibri run, noProof). The fix removes that
# character's run entirely while leaving the (now empty) paragraph -
# and its paragraph formatting / mark - in place.
$p1 = $d.Paragraphs(1)
$runRange = $d.Range($p1.Range.Start, $p1.Range.End - 1)

if ($runRange.Text -eq "H") {
    $runRange.Delete()
} else {
    # Defensive fallback in case the stray run has moved: locate the
    # lone "H" run directly via Find and remove just its text.
    $found = $d.Content
    $found.Find.Execute("H", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)
}
